$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.776.88"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.453.43"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'581.12"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "'148.56"
$ws.Range("E6").Value = "  +9.43%  "
$ws.Range("D7").Value = "3.451.47"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "4.046.58"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "'28.10"
$ws.Range("E14").Value = "  +8.05%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "3.480.97"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "61.874.89"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  +8.79%  "
$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "'9.45"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'385.78"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").Value = "3.600.33"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "'72.66"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'0.180"
$ws.Range("E29").Value = "  +9.13%  "
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'1.51"
$ws.Range("E32").Value = "  -13.63%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D36").Value = "'23.97"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").Value = "'7.08"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").Value = "'5.23"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "'166.02"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'0.0790"
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("D42").Value = "'26.14"
$ws.Range("E42").Value = "  +9.78%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'42.38"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").Value = "'1.72"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "2.594.80"
$ws.Range("E49").Value = "  +10.48%  "
$ws.Range("D50").Value = "'6.97"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "'23.32"
$ws.Range("E51").Value = "  +0.70%  "
